$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number must be forced to text
# (matching the source workbook, which stores these as inline strings)
# so Excel does not silently convert them to the Number type.
$textForceCells = @("D5", "D6", "D8", "D9", "D14", "D19", "D20", "D21", "D22", "D24", "D26", "D28", "D29", "D31", "D33", "D34", "D35", "D36", "D37", "D40", "D41", "D43", "D47", "D48", "D49", "D51")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "57.269.73"
$ws.Range("D3").Value = "3.065.56"
$ws.Range("E3").Value = "  +0.31%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "515.30"
$ws.Range("E5").Value = "  -0.32%  "
$ws.Range("D6").Value = "140.76"
$ws.Range("E6").Value = "  -0.38%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "0.432"
$ws.Range("E8").Value = "  -1.86%  "
$ws.Range("D9").Value = "7.25"
$ws.Range("E9").Value = "  +0.12%  "
$ws.Range("E10").Value = "  -1.68%  "
$ws.Range("E11").Value = "  -1.92%  "
$ws.Range("D12").Value = "3.597.13"
$ws.Range("E12").Value = "  +0.48%  "
$ws.Range("E13").Value = "  +2.54%  "
$ws.Range("D14").Value = "25.50"
$ws.Range("E14").Value = "  -5.28%  "
$ws.Range("E15").Value = "  -2.50%  "
$ws.Range("D16").Value = "57.388.09"
$ws.Range("E16").Value = "  +0.29%  "
$ws.Range("D17").Value = "3.070.97"
$ws.Range("E17").Value = "  +0.45%  "
$ws.Range("E18").Value = "  -1.51%  "
$ws.Range("D19").Value = "13.04"
$ws.Range("E19").Value = "  -3.02%  "
$ws.Range("D20").Value = "8.08"
$ws.Range("E20").Value = "  -0.44%  "
$ws.Range("D21").Value = "332.47"
$ws.Range("E21").Value = "  -1.11%  "
$ws.Range("D22").Value = "0.996"
$ws.Range("E22").Value = "  -0.48%  "
$ws.Range("E23").Value = "  -1.72%  "
$ws.Range("D24").Value = "65.65"
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("E25").Value = "  +2.83%  "
$ws.Range("D26").Value = "0.995"
$ws.Range("E26").Value = "  -0.85%  "
$ws.Range("D27").Value = "0.0₃0906"
$ws.Range("E27").Value = "  +0.11%  "
$ws.Range("D28").Value = "6.33"
$ws.Range("E28").Value = "  -5.97%  "
$ws.Range("D29").Value = "7.14"
$ws.Range("E29").Value = "  -0.83%  "
$ws.Range("E30").Value = "  -0.12%  "
$ws.Range("D31").Value = "20.78"
$ws.Range("E31").Value = "  -0.13%  "
$ws.Range("E32").Value = "  -3.11%  "
$ws.Range("D33").Value = "154.82"
$ws.Range("E33").Value = "  +1.18%  "
$ws.Range("D34").Value = "27.31"
$ws.Range("E34").Value = "  +8.01%  "
$ws.Range("D35").Value = "4.48"
$ws.Range("E35").Value = "  -5.42%  "
$ws.Range("D36").Value = "5.83"
$ws.Range("E36").Value = "  -2.09%  "
$ws.Range("D37").Value = "1.27"
$ws.Range("E37").Value = "  -0.97%  "
$ws.Range("E38").Value = "  -0.51%  "
$ws.Range("D39").Value = "3.106.92"
$ws.Range("E39").Value = "  +0.51%  "
$ws.Range("D40").Value = "36.80"
$ws.Range("E40").Value = "  -0.92%  "
$ws.Range("D41").Value = "3.87"
$ws.Range("E41").Value = "  -0.66%  "
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("D43").Value = "0.658"
$ws.Range("E43").Value = "  -1.16%  "
$ws.Range("D44").Value = "2.271.23"
$ws.Range("E44").Value = "  +2.43%  "
$ws.Range("E45").Value = "  +6.52%  "
$ws.Range("E46").Value = "  -1.99%  "
$ws.Range("B47").Value = "Cosmos"
$ws.Range("C47").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D47").Value = "5.87"
$ws.Range("E47").Value = "  -2.73%  "
$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D48").Value = "0.928"
$ws.Range("E48").Value = "  -4.07%  "
$ws.Range("D49").Value = "19.80"
$ws.Range("E49").Value = "  -2.64%  "
$ws.Range("E50").Value = "  +0.21%  "
$ws.Range("D51").Value = "248.37"
$ws.Range("E51").Value = "  +5.03%  "

foreach ($addr in $textForceCells) {
    $ws.Range($addr).Style = "Normal"
}
